$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalized cell-count values for B1:G23, per Grant's request
$ws.Range("B1").Value = 1417.8624717407686
$ws.Range("C1").Value = 1585.5476391740131
$ws.Range("D1").Value = 1476.573705179283
$ws.Range("E1").Value = 1734.4564315352698
$ws.Range("F1").Value = 1676.4533046884869
$ws.Range("G1").Value = 1491.80699023631
$ws.Range("B2").Value = 1527.1088922381311
$ws.Range("C2").Value = 1726.1219659461419
$ws.Range("D2").Value = 1621.9316283254082
$ws.Range("E2").Value = 1785.3202565069785
$ws.Range("F2").Value = 1710.3963098698341
$ws.Range("G2").Value = 1512.5265317673698
$ws.Range("B3").Value = 1563.7844762622456
$ws.Range("C3").Value = 1788.2362033570828
$ws.Range("D3").Value = 1693.3787430921477
$ws.Range("E3").Value = 1817.5340123223941
$ws.Range("F3").Value = 1710.3963098698341
$ws.Range("G3").Value = 1518.3280033960666
$ws.Range("B4").Value = 1643.378296910324
$ws.Range("C4").Value = 1878.1383890834443
$ws.Range("D4").Value = 1726.2279912607637
$ws.Range("E4").Value = 1868.3978372941028
$ws.Range("F4").Value = 1789.8726146846957
$ws.Range("G4").Value = 1582.144191311731
$ws.Range("B5").Value = 1701.9031650339109
$ws.Range("C5").Value = 1950.8774302620459
$ws.Range("D5").Value = 1808.3511116823033
$ws.Range("E5").Value = 1920.9571230982019
$ws.Range("F5").Value = 1888.3901175281178
$ws.Range("G5").Value = 1652.5906325173344
$ws.Range("B6").Value = 1797.8839487565938
$ws.Range("C6").Value = 2046.500664170994
$ws.Range("D6").Value = 1868.3009895900273
$ws.Range("E6").Value = 1988.7755563938138
$ws.Range("F6").Value = 1978.6288386199919
$ws.Range("G6").Value = 1726.3522003679075
$ws.Range("B7").Value = 1874.3564431047475
$ws.Range("C7").Value = 2134.7682647023307
$ws.Range("D7").Value = 1924.1447114766743
$ws.Range("E7").Value = 2087.9600150886458
$ws.Range("F7").Value = 2025.8178946038161
$ws.Range("G7").Value = 1819.1757464270556
$ws.Range("B8").Value = 1961.7535795026374
$ws.Range("C8").Value = 2194.4306243207343
$ws.Range("D8").Value = 1977.5247397506751
$ws.Range("E8").Value = 2154.0829875518671
$ws.Range("F8").Value = 2113.5729811702258
$ws.Range("G8").Value = 1902.8826942125374
$ws.Range("B9").Value = 2045.2490580256217
$ws.Range("C9").Value = 2272.0734210844103
$ws.Range("D9").Value = 2062.932784989076
$ws.Range("E9").Value = 2252.4197158305042
$ws.Range("F9").Value = 2195.5329205105518
$ws.Range("G9").Value = 1967.5276637894442
$ws.Range("B10").Value = 2113.1379050489827
$ws.Range("C10").Value = 2392.2154329187297
$ws.Range("D10").Value = 2177.9051535792319
$ws.Range("E10").Value = 2370.2542436816298
$ws.Range("F10").Value = 2282.4601289018065
$ws.Range("G10").Value = 2054.5497382198955
$ws.Range("B11").Value = 2324.6077618688769
$ws.Range("C11").Value = 2541.7799782634952
$ws.Range("D11").Value = 2336.4027759928035
$ws.Range("E11").Value = 2560.9935873255376
$ws.Range("F11").Value = 2487.7739163401993
$ws.Range("G11").Value = 2251.7997735955855
$ws.Range("B12").Value = 2456.4837980406933
$ws.Range("C12").Value = 2728.1226904963173
$ws.Range("D12").Value = 2472.727155892559
$ws.Range("E12").Value = 2745.7988180560797
$ws.Range("F12").Value = 2649.2101604953868
$ws.Range("G12").Value = 2375.2882411207024
$ws.Range("B13").Value = 2630.4977392614919
$ws.Range("C13").Value = 2885.0428692186938
$ws.Range("D13").Value = 2669.8226449042545
$ws.Range("E13").Value = 2928.9085879542313
$ws.Range("F13").Value = 2846.2451661822311
$ws.Range("G13").Value = 2526.1265034668181
$ws.Range("B14").Value = 2788.1247174076866
$ws.Range("C14").Value = 3086.0968482067387
$ws.Range("D14").Value = 2808.6107184166563
$ws.Range("E14").Value = 3088.2819061989189
$ws.Range("F14").Value = 2991.9517250094773
$ws.Range("G14").Value = 2686.9101457478428
$ws.Range("B15").Value = 2909.8564431047475
$ws.Range("C15").Value = 3266.718512256974
$ws.Range("D15").Value = 2990.1028145482592
$ws.Range("E15").Value = 3238.3301898654595
$ws.Range("F15").Value = 3140.1419183621883
$ws.Range("G15").Value = 2882.5026178010476
$ws.Range("B16").Value = 2943.4107008289375
$ws.Range("C16").Value = 3263.4493418669244
$ws.Range("D16").Value = 3077.9745533993064
$ws.Range("E16").Value = 3312.0827360744374
$ws.Range("F16").Value = 3277.5696954378864
$ws.Range("G16").Value = 2836.9196264327161
$ws.Range("B17").Value = 3030.0275056518462
$ws.Range("C17").Value = 3323.1117014853285
$ws.Range("D17").Value = 3100.9690271173376
$ws.Range("E17").Value = 3331.5805356469259
$ws.Range("F17").Value = 3241.1430557310746
$ws.Range("G17").Value = 2889.9616527522294
$ws.Range("B18").Value = 3024.565184626978
$ws.Range("C18").Value = 3334.5537978505017
$ws.Range("D18").Value = 3118.2148824058609
$ws.Range("E18").Value = 3335.8191877279014
$ws.Range("F18").Value = 3271.7745482118025
$ws.Range("G18").Value = 2899.0782510258955
$ws.Range("B19").Value = 3062.0211002260739
$ws.Range("C19").Value = 3367.2455017509965
$ws.Range("D19").Value = 3140.3881249196766
$ws.Range("E19").Value = 3347.6874135546336
$ws.Range("F19").Value = 3322.2751168962459
$ws.Range("G19").Value = 2925.5992641856524
$ws.Range("B20").Value = 3057.3391107761868
$ws.Range("C20").Value = 3361.5244535684101
$ws.Range("D20").Value = 3154.3490553913384
$ws.Range("E20").Value = 3348.5351439708288
$ws.Range("F20").Value = 3279.2254517881961
$ws.Range("G20").Value = 2947.9763690391969
$ws.Range("B21").Value = 3058.1194423511679
$ws.Range("C21").Value = 3376.2357203236329
$ws.Range("D21").Value = 3141.2093561238921
$ws.Range("E21").Value = 3321.4077706525841
$ws.Range("F21").Value = 3281.7090863136605
$ws.Range("G21").Value = 2943.0036790717427
$ws.Range("B22").Value = 3057.3391107761868
$ws.Range("C22").Value = 3392.5815722738803
$ws.Range("D22").Value = 3181.4496851304466
$ws.Range("E22").Value = 3369.7284043757072
$ws.Range("F22").Value = 3274.2581827372669
$ws.Range("G22").Value = 2924.77048252441
$ws.Range("B23").Value = 3090.1130369253956
$ws.Range("C23").Value = 3386.8605240912939
$ws.Range("D23").Value = 3190.4832283768155
$ws.Range("E23").Value = 3370.5761347919024
$ws.Range("F23").Value = 3275.086060912422
$ws.Range("G23").Value = 2927.2568275081371

# Restore the selection that Excel records after this edit (B1:G23, active cell B1)
$ws.Range("B1:G23").Select()
